$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add a new "HN vs LN" contrast row (row 23) to each of the three result
# sheets: MaxPixelCount, InflectionPoint, MaxGrowthRate.
# ---------------------------------------------------------------------------

$ws1 = $wb.Worksheets.Item("MaxPixelCount")
$ws1.Cells.Item(23, 1).Value = "HN vs LN"
$ws1.Cells.Item(23, 2).Value = 1213488.53350966
$ws1.Cells.Item(23, 3).Value = 143549.54044249601
$ws1.Cells.Item(23, 4).Value = 24
$ws1.Cells.Item(23, 5).Value = 8.4534477071054308
$ws1.Cells.Item(23, 6).Value = 0.0000000117320274406411
$ws1.Range("F23").Interior.Color = 65535
$ws1.Range("F23").NumberFormat = "0.00E+00"

$ws2 = $wb.Worksheets.Item("InflectionPoint")
$ws2.Cells.Item(23, 1).Value = "HN vs LN"
$ws2.Cells.Item(23, 2).Value = 29.390894492466401
$ws2.Cells.Item(23, 3).Value = 5.09660742963404
$ws2.Cells.Item(23, 4).Value = 24
$ws2.Cells.Item(23, 5).Value = 5.7667565921546302
$ws2.Cells.Item(23, 6).Value = 0.0000060624960686767401
$ws2.Range("F23").Interior.Color = 65535
$ws2.Range("F23").NumberFormat = "0.00E+00"

$ws3 = $wb.Worksheets.Item("MaxGrowthRate")
$ws3.Cells.Item(23, 1).Value = "HN vs LN"
$ws3.Cells.Item(23, 2).Value = 27925.204634959598
$ws3.Cells.Item(23, 3).Value = 3252.5436599488498
$ws3.Cells.Item(23, 4).Value = 24
$ws3.Cells.Item(23, 5).Value = 8.5856509718300806
$ws3.Cells.Item(23, 6).Value = 0.0000000088371146908655803
$ws3.Range("F23").Interior.Color = 65535
$ws3.Range("F23").NumberFormat = "0.00E+00"

# Column A on MaxGrowthRate gains an explicit custom width (it had none
# before); 17.1 characters renders as the stored width of 18.
$ws3.Columns.Item(1).ColumnWidth = 17.1

# ---------------------------------------------------------------------------
# Selection / active-sheet bookkeeping: every sheet ends up with F23
# selected, and the active tab moves from MaxGrowthRate to InflectionPoint.
# ---------------------------------------------------------------------------

$ws1.Range("F23").Select() | Out-Null
$ws3.Range("F23").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("F23").Select() | Out-Null
